$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the title
#    (Heading1) paragraph at the top of the document.
# ------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item(2)
$metaPara.Style = "Normal"

$metaLabel = "Meta description"
$metaRest  = ": Read our review of Alchemistress slot, including bonus features and free spins. Play this Harry-Potter inspired game for free now."
$metaFull  = $metaLabel + $metaRest

$insertStart = $metaPara.Range.Start
$insertionRange = $d.Range($insertStart, $insertStart)
$insertionRange.Text = $metaFull

$boldRange = $d.Range($insertStart, $insertStart + $metaLabel.Length)
$boldRange.Bold = 1

# ------------------------------------------------------------------
# 2) Near the end of the document, remove the duplicated bold
#    "Play Alchemistress Slot | Free Review & Free Spins" paragraph,
#    then update the italic paragraph's text with the new image
#    generation prompt.
# ------------------------------------------------------------------
$count = $d.Paragraphs.Count
$boldDupPara = $d.Paragraphs.Item($count - 1)
$boldDupPara.Range.Delete()

$count = $d.Paragraphs.Count
$italicPara = $d.Paragraphs.Item($count)
$italicStart = $italicPara.Range.Start
$italicEnd = $italicPara.Range.End
$italicRange = $d.Range($italicStart, $italicEnd)

$newPrompt = 'Create a feature image for the game "Alchemistress" that fits the theme of the game. The image should be in a cartoon style and feature a happy Maya warrior with glasses. The Maya warrior should be standing in a candlelit room with a game board and potion vials in the background. The warrior should be holding a wand and wearing a wizard''s hat. The image should be bright and colorful to capture the excitement of the game. The text "Alchemistress" should be prominently displayed in a fun and playful font.'

$italicRange.Text = $newPrompt

Write-Output "done"
